$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2130.86
$ws.Range("D2").Value = 2939.46
$ws.Range("E2").Value = 2932.53

$ws.Range("C3").Value = 207.81
$ws.Range("D3").Value = 512.86
$ws.Range("E3").Value = 1612.8915

$ws.Range("C4").Value = 102.86
$ws.Range("D4").Value = 98.14
$ws.Range("E4").Value = 293.253

$ws.Range("C5").Value = 101.98
$ws.Range("D5").Value = 64.34999999999999
$ws.Range("E5").Value = 439.8795

$ws.Range("C6").Value = 14.08
$ws.Range("D6").Value = 42.25

$ws.Range("C7").Value = 373.54
$ws.Range("D7").Value = 272.53

$ws.Range("C8").Value = 776.75
$ws.Range("D8").Value = 858.54

$ws.Range("C9").Value = 2.43
$ws.Range("D9").Value = 2.15
$ws.Range("E9").Value = 29.3253

$ws.Range("C10").Value = 41.02
$ws.Range("D10").Value = 19.54
$ws.Range("E10").Value = 293.253

$ws.Range("C11").Value = 21.37
$ws.Range("D11").Value = 20.08
$ws.Range("E11").Value = 175.9518

$ws.Range("C12").Value = 3454.54
$ws.Range("D12").Value = 2196.35
$ws.Range("E12").Value = 2000

$ws.Range("C13").Value = 1879.42
$ws.Range("D13").Value = 4273.16

$ws.Range("C14").Value = 15
$ws.Range("D14").Value = 16.57

$ws.Range("C15").Value = 172.62
$ws.Range("D15").Value = 374.65

$ws.Range("C16").Value = 0.6899999999999999
$ws.Range("D16").Value = 1.26

$ws.Range("C17").Value = 1.1
$ws.Range("D17").Value = 1.44

$ws.Range("C18").Value = 0.6899999999999999
$ws.Range("D18").Value = 1.29

$ws.Range("C19").Value = 9.23
$ws.Range("D19").Value = 10.81

$ws.Range("C20").Value = 5.93
$ws.Range("D20").Value = 5.35

$ws.Range("C21").Value = 16.92
$ws.Range("D21").Value = 145.34

$ws.Range("C22").Value = 236.89
$ws.Range("D22").Value = 647.35

$ws.Range("C23").Value = 0.9
$ws.Range("D23").Value = 1.48

$ws.Range("C24").Value = 356.56
$ws.Range("D24").Value = 568.53

$ws.Range("C25").Value = 1262.93
$ws.Range("D25").Value = 1416.36

$ws.Range("C26").Value = 18.83
$ws.Range("D26").Value = 19.15
